# "R12 Superheated Vapour Data"
#
# The "Condenser Temp [C]" column (column C) is no longer collected /
# needed, so the whole column is removed from the data table on Sheet1
# (Excel shifts every later column one position to the left, updating
# the used range, column-width definitions and all cell references
# automatically).
#
# The column that used to be labelled "Temperature CompressorOutlet 3"
# (old column Q, now column P after the deletion above) actually holds
# the Condenser Outlet temperature for the superheated-vapour run, so
# its header is corrected/renamed to "Temperature Condenser Outlet 3".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the entire "Condenser Temp [C]" column (old column C).
$ws.Columns("C").Delete()

# Rename the header that is now in column P.
$ws.Range("P1").Value = "Temperature Condenser Outlet 3"

# Restore the cursor/selection the author left the sheet with.
$ws.Range("S1").Select()
